# Updates cryptos list: Price/Volume(1h) refresh + two row swaps (rows 25/26, 49/50)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.359.38"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "2.322.72"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.43%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("E9").Value = "  -3.22%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "2.738.15"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "56.343.41"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "2.324.41"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "325.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.89%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.162"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("E27").Value = "  +2.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.03%  "
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("D30").Value = "0.0₃0716"
$ws.Range("E30").Value = "  -4.34%  "
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("E36").Value = "  -3.87%  "
$ws.Range("E37").Value = "  -2.78%  "
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "148.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.372"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "276.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0926"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0213"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.374"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.30%  "
